$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The KiCad BOM export was regenerated with a new "References" column
# (the full list of reference designators grouped into each BOM line),
# inserted right after the "Count" column (C) and before "Name" (old D).
# Inserting a whole column shifts Name/Footprint/Description (old D/E/F)
# and the last column (old H) one slot to the right.
$ws.Columns("D").Insert()

# New column gets the same text formatting as column A (Reference).
$ws.Range("A1:A13").Copy() | Out-Null
$ws.Range("D1:D13").PasteSpecial(-4122) | Out-Null   ; # xlPasteFormats
$excel.CutCopyMode = $false

# Header + per-row reference lists. Rows 3-13 duplicate column A; row 2's
# grouping was edited (C16 dropped, C19 added) and row 1 is the header.
$ws.Range("D1").Value = "References"
$ws.Range("D2").Value = "C2, C3, C4, C5, C6, C7, C8, C9, C12, C13, C14, C15, C17, C19"
$ws.Range("D3").Value = "C10, C11, C18"
$ws.Range("D4").Value = "C17"
$ws.Range("D5").Value = "C1"
$ws.Range("D6").Value = "D1"
$ws.Range("D7").Value = "J2"
$ws.Range("D8").Value = "R1"
$ws.Range("D9").Value = "U1"
$ws.Range("D10").Value = "U2"
$ws.Range("D11").Value = "U3"
$ws.Range("D12").Value = "U4"
$ws.Range("D13").Value = "U5"

# Column widths: "Count" (C) narrows, the new "References" column (D)
# matches "Reference" (A), and "Name"/"Footprint" (now E/F) both end up
# at the old "Footprint" best-fit width.
$ws.Columns("C").ColumnWidth = 6.5
$ws.Columns("D").ColumnWidth = 50.67
$ws.Columns("E").ColumnWidth = 35.5
$ws.Columns("F").ColumnWidth = 35.5

# The trailing blank formatted cells (old H10/I10) shift to I10/J10; drop
# the redundant one so only I10 remains blank, same as the original file
# only ever carried one spare blank cell at the end of that row.
$ws.Range("J10").Clear() | Out-Null

# Selection moves to D14 (just below the newly filled column) after the edit.
$ws.Range("D14").Select()
